# Apply updated Betfair Back/Lay odds values for rows 2-15 (F:AO)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.67
$ws.Range("G2").Value = 1.82
$ws.Range("H2").Value = 6.4
$ws.Range("I2").Value = 8.199999999999999
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 3.75
$ws.Range("L2").Value = 1.57
$ws.Range("M2").Value = 1.12
$ws.Range("N2").Value = 2.46
$ws.Range("O2").Value = 1.56
$ws.Range("P2").Value = 1.49
$ws.Range("Q2").Value = 2.64
$ws.Range("R2").Value = 1.17
$ws.Range("S2").Value = 5.7
$ws.Range("T2").Value = 2.42
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 1.14
$ws.Range("W2").Value = 2.2
$ws.Range("X2").Value = 8.800000000000001
$ws.Range("Y2").Value = 18.5
$ws.Range("AB2").Value = 6.4
$ws.Range("AC2").Value = 9.800000000000001
$ws.Range("AD2").Value = 38
$ws.Range("AF2").Value = 10
$ws.Range("AG2").Value = 12
$ws.Range("AH2").Value = 36
$ws.Range("AJ2").Value = 23
$ws.Range("AK2").Value = 28
$ws.Range("AN2").Value = 22
# Row 3
$ws.Range("G3").Value = 3.45
$ws.Range("H3").Value = 2.2
$ws.Range("I3").Value = 2.48
$ws.Range("J3").Value = 3.5
$ws.Range("L3").Value = 1.28
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 2.12
$ws.Range("Q3").Value = 1.76
$ws.Range("S3").Value = 2.84
$ws.Range("T3").Value = 1.63
$ws.Range("U3").Value = 2.28
$ws.Range("V3").Value = 1.69
$ws.Range("W3").Value = 1.38
$ws.Range("X3").Value = 22
$ws.Range("Y3").Value = 14
$ws.Range("AA3").Value = 36
$ws.Range("AB3").Value = 15.5
$ws.Range("AE3").Value = 25
$ws.Range("AF3").Value = 28
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 36
$ws.Range("AJ3").Value = 60
$ws.Range("AK3").Value = 40
$ws.Range("AL3").Value = 1000
$ws.Range("AO3").Value = 18.5
# Row 5
$ws.Range("G5").Value = 1.77
$ws.Range("L5").Value = 1.5
$ws.Range("M5").Value = 1.11
$ws.Range("O5").Value = 1.57
$ws.Range("Y5").Value = 19.5
$ws.Range("AC5").Value = 10
$ws.Range("AD5").Value = 36
$ws.Range("AG5").Value = 13.5
$ws.Range("AH5").Value = 40
$ws.Range("AJ5").Value = 21
$ws.Range("AK5").Value = 30
$ws.Range("AN5").Value = 21
# Row 6
$ws.Range("J6").Value = 3.8
$ws.Range("K6").Value = 3.95
$ws.Range("N6").Value = 3.4
$ws.Range("O6").Value = 1.38
$ws.Range("P6").Value = 1.83
$ws.Range("Q6").Value = 2.1
$ws.Range("R6").Value = 1.31
$ws.Range("S6").Value = 3.9
$ws.Range("T6").Value = 1.94
$ws.Range("U6").Value = 1.93
$ws.Range("AL6").Value = 85
# Row 7
$ws.Range("F7").Value = 1.37
$ws.Range("G7").Value = 1.5
$ws.Range("H7").Value = 7.6
$ws.Range("J7").Value = 3.65
$ws.Range("K7").Value = 7.2
$ws.Range("N7").Value = 2.06
$ws.Range("P7").Value = 1.55
$ws.Range("Q7").Value = 2.02
$ws.Range("S7").Value = 2.92
$ws.Range("W7").Value = 3
# Row 8
$ws.Range("F8").Value = 1.5
$ws.Range("H8").Value = 1.09
$ws.Range("I8").Value = 10
$ws.Range("J8").Value = 3.3
$ws.Range("K8").Value = 4.8
$ws.Range("Q8").Value = 1.94
# Row 9
$ws.Range("G9").Value = 1.55
$ws.Range("H9").Value = 6.6
$ws.Range("J9").Value = 4.8
$ws.Range("N9").Value = 5
$ws.Range("O9").Value = 1.23
$ws.Range("P9").Value = 2.36
$ws.Range("Q9").Value = 1.68
$ws.Range("R9").Value = 1.53
$ws.Range("S9").Value = 2.72
$ws.Range("T9").Value = 1.8
$ws.Range("W9").Value = 2.8
$ws.Range("Y9").Value = 980
$ws.Range("AC9").Value = 10.5
$ws.Range("AI9").Value = 85
$ws.Range("AJ9").Value = 14.5
$ws.Range("AN9").Value = 7
$ws.Range("AO9").Value = 95
# Row 10
$ws.Range("F10").Value = 1.68
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 7.6
$ws.Range("J10").Value = 2.94
$ws.Range("N10").Value = 1.91
$ws.Range("P10").Value = 1.46
$ws.Range("Q10").Value = 2.16
$ws.Range("S10").Value = 2.16
$ws.Range("V10").Value = 1.15
$ws.Range("W10").Value = 1.87
# Row 11
$ws.Range("F11").Value = 2.46
$ws.Range("G11").Value = 2.54
$ws.Range("H11").Value = 3.55
$ws.Range("M11").Value = 1.13
$ws.Range("N11").Value = 2.72
$ws.Range("O11").Value = 1.54
$ws.Range("Q11").Value = 2.64
$ws.Range("T11").Value = 2.12
$ws.Range("U11").Value = 1.82
$ws.Range("V11").Value = 1.37
$ws.Range("W11").Value = 1.65
$ws.Range("X11").Value = 8.800000000000001
$ws.Range("Z11").Value = 22
$ws.Range("AC11").Value = 6.8
$ws.Range("AD11").Value = 16
$ws.Range("AF11").Value = 13.5
$ws.Range("AH11").Value = 27
$ws.Range("AI11").Value = 85
$ws.Range("AJ11").Value = 36
$ws.Range("AK11").Value = 36
$ws.Range("AL11").Value = 65
$ws.Range("AM11").Value = 180
$ws.Range("AN11").Value = 38
# Row 12
$ws.Range("F12").Value = 1.94
$ws.Range("G12").Value = 2.16
$ws.Range("I12").Value = 5.6
$ws.Range("J12").Value = 3.25
$ws.Range("Q12").Value = 2.02
# Row 13
$ws.Range("G13").Value = 2.34
$ws.Range("H13").Value = 3.8
$ws.Range("I13").Value = 4.1
$ws.Range("P13").Value = 1.56
$ws.Range("Q13").Value = 2.52
$ws.Range("V13").Value = 1.32
$ws.Range("W13").Value = 1.75
$ws.Range("Z13").Value = 980
$ws.Range("AD13").Value = 980
$ws.Range("AF13").Value = 980
$ws.Range("AG13").Value = 980
$ws.Range("AH13").Value = 980
$ws.Range("AJ13").Value = 980
$ws.Range("AK13").Value = 980
# Row 14
$ws.Range("F14").Value = 1.7
$ws.Range("I14").Value = 6.4
$ws.Range("J14").Value = 3.75
$ws.Range("K14").Value = 4.4
$ws.Range("P14").Value = 1.88
$ws.Range("Q14").Value = 1.98
$ws.Range("T14").Value = 1.92
$ws.Range("X14").Value = 980
$ws.Range("Y14").Value = 980
$ws.Range("AB14").Value = 980
$ws.Range("AC14").Value = 11
$ws.Range("AD14").Value = 980
$ws.Range("AF14").Value = 980
$ws.Range("AG14").Value = 980
$ws.Range("AH14").Value = 980
$ws.Range("AJ14").Value = 980
$ws.Range("AK14").Value = 980
$ws.Range("AL14").Value = 980
$ws.Range("AN14").Value = 980
# Row 15
$ws.Range("Q15").Value = 1.78
$ws.Range("R15").Value = 1.44
$ws.Range("S15").Value = 2.92
$ws.Range("T15").Value = 1.68
$ws.Range("U15").Value = 2.22
$ws.Range("Z15").Value = 980
$ws.Range("AA15").Value = 75
$ws.Range("AE15").Value = 44
$ws.Range("AI15").Value = 980
$ws.Range("AJ15").Value = 980
$ws.Range("AL15").Value = 980
$ws.Range("AM15").Value = 85
$ws.Range("AN15").Value = 13.5
$ws.Range("AO15").Value = 1000
